$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row 3, column F label (#说明 -> 套装ID) ---
$ws.Range("F3").Value = "套装ID"

# --- Row 4: insert new SuitId field name in column F ---
$ws.Range("F4").Value = "SuitId"

# --- Row 5: insert new "int" type marker in column F ---
$ws.Range("F5").Value = "int"

# --- Rows 6-10: add SuitId value (20011) in column F ---
$ws.Range("F6").Value = 20011
$ws.Range("F7").Value = 20011
$ws.Range("F8").Value = 20011
$ws.Range("F9").Value = 20011
$ws.Range("F10").Value = 20011

# --- Fully remove the now-unused, style-only L6 cell ---
$ws.Range("L6").Clear()

# --- Prepare new rows 11-15 by copying formatting from row 6 (skip column L, which stays blank/unset) ---
$ws.Range("C6:K6").Copy()
$ws.Range("C11:K15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M6:P6").Copy()
$ws.Range("M11:P15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate new rows 11-15 (Suit 20012 variants) ---
$ws.Range("C11").Value = 20021
$ws.Range("D11").Value = 2001
$ws.Range("E11").Value = "火球·狂暴"
$ws.Range("F11").Value = 20012
$ws.Range("G11").Value = "增加火球20%伤害系数，无上限"
$ws.Range("H11").Value = 99
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 20
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0

$ws.Range("C12").Value = 20022
$ws.Range("D12").Value = 2001
$ws.Range("E12").Value = "火球·多重"
$ws.Range("F12").Value = 20012
$ws.Range("G12").Value = "增加火球术一个攻击目标，最多生效3个"
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 1
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0

$ws.Range("C13").Value = 20023
$ws.Range("D13").Value = 2001
$ws.Range("E13").Value = "火球·急速"
$ws.Range("F13").Value = 20012
$ws.Range("G13").Value = "增加火球术1秒CD时间，最多生效1个"
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0

$ws.Range("C14").Value = 20024
$ws.Range("D14").Value = 2001
$ws.Range("E14").Value = "火球·灼烧"
$ws.Range("F14").Value = 20012
$ws.Range("G14").Value = "火球术增加烧伤伤害，灼烧系数为80%,最多叠加3重"
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 101

$ws.Range("C15").Value = 20025
$ws.Range("D15").Value = 2001
$ws.Range("E15").Value = "火球·固伤"
$ws.Range("F15").Value = 20012
$ws.Range("G15").Value = "增加火球术100点固定伤害，无上限"
$ws.Range("H15").Value = 99
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 100
$ws.Range("P15").Value = 0

# --- Update selection to mimic final author cursor position (cosmetic) ---
$ws.Range("E23").Select()
